{"js": "// The document contains a link to the course's GitHub repository whose\n// name is suffixed with the year, e.g.:\n//   https://github.com/wesarmour/CWM-in-HPC-and-Scientific-Computing-2022.git\n// The practical notes were updated for this year's course, so the \"2022\"\n// needs to become \"2023\" (the rest of the paragraph, including the \".git\"\n// suffix and the surrounding \"git clone \" text, stays the same).\nconst results = context.document.body.search(\"CWM-in-HPC-and-Scientific-Computing-2022\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\n    \"CWM-in-HPC-and-Scientific-Computing-2023\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# The practical notes' git-clone instructions point at a GitHub repo whose\n# name is suffixed with the year, e.g.:\n#   https://github.com/wesarmour/CWM-in-HPC-and-Scientific-Computing-2022.git\n# This year's update simply bumps that year from 2022 to 2023 (everything\n# else in the sentence, including the \".git\" suffix, is unchanged).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"CWM-in-HPC-and-Scientific-Computing-2022\"\n$find.Replacement.Text = \"CWM-in-HPC-and-Scientific-Computing-2023\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.Execute($null, $false, $false, $false, $null, $null, $true, $null, $null, $null, 2) | Out-Null\n"}
